# Add two new columns (I: "I0", J: "IF") to the table on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing last header cell (H1, bold/border/centered
# style) onto the two new header cells so they match the rest of the header row.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

# Header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for the two columns
$ws.Range("I2").Value = 4
$ws.Range("J2").Value = 6
$ws.Range("I3").Value = 6
$ws.Range("J3").Value = 7
